# Add a new "topup" column (Q) with a sample value, matching the
# "Tambah Pelanggan ... Top Up Pelanggan" test-data commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell + sample data cell in column Q, mirroring the existing
# header/data row layout (A1:P1 headers, A2:P2 sample row).
$ws.Range("Q1").Value = "topup"
$ws.Range("Q2").Value = "10000"

# Reflect the view state recorded in the saved workbook: the new cell is
# selected and the sheet has been scrolled right so column Q is visible.
$ws.Range("Q2").Select()
$excel.ActiveWindow.ScrollColumn = 10
